$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 1263.9166
$ws.Range("J5").Value = 613.3333
$ws.Range("L5").Value = 613.3333
$ws.Range("N5").Value = -843.3333
$ws.Range("H16").Value = 25666.334
$ws.Range("J16").Value = 37999.5
$ws.Range("L16").Value = 37999.5
$ws.Range("N16").Value = -38459.5
$ws.Range("H48").Value = 3880
$ws.Range("J48").Value = 4000
$ws.Range("L48").Value = 12000
$ws.Range("N48").Value = -12584
$ws.Range("H56").Value = 3880
$ws.Range("J56").Value = 4000
$ws.Range("L56").Value = 12000
$ws.Range("N56").Value = -13068
$ws.Range("H92").Value = 903.4286
$ws.Range("I92").Value = 1048.6666
$ws.Range("J92").Value = 642
$ws.Range("K92").Value = 1048.6666
$ws.Range("L92").Value = 642
$ws.Range("M92").Value = 199.3334
$ws.Range("N92").Value = -3138
$ws.Range("H94").Value = 4443.5
$ws.Range("I94").Value = 4443.5
$ws.Range("K94").Value = 4443.5
$ws.Range("M94").Value = -3992.5
$ws.Range("H96").Value = 62507292
$ws.Range("I96").Value = 5934.8335
$ws.Range("J96").Value = 100008104
$ws.Range("K96").Value = 17804.5005
$ws.Range("L96").Value = 300024312
$ws.Range("M96").Value = -16431.5005
$ws.Range("N96").Value = -300027058
$ws.Range("H99").Value = 1939.5
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("H100").Value = 1384.8823
$ws.Range("I100").Value = 1701.1
$ws.Range("J100").Value = 933.1429000000001
$ws.Range("K100").Value = 1701.1
$ws.Range("L100").Value = 933.1429000000001
$ws.Range("M100").Value = -1160.1
$ws.Range("N100").Value = -2015.1429
$ws.Range("H101").Value = 692.06665
$ws.Range("J101").Value = 957.5
$ws.Range("L101").Value = 2872.5
$ws.Range("N101").Value = -6116.5
$ws.Range("H103").Value = 4459.5
$ws.Range("J103").Value = 4459
$ws.Range("L103").Value = 13377
$ws.Range("N103").Value = -14549
$ws.Range("H116").Value = 5962.913
$ws.Range("I116").Value = 3027.7778
$ws.Range("J116").Value = 7849.7856
$ws.Range("K116").Value = 3027.7778
$ws.Range("L116").Value = 7849.7856
$ws.Range("M116").Value = 414.2222000000002
$ws.Range("N116").Value = -14733.7856
$ws.Range("H123").Value = 34044.668
$ws.Range("J123").Value = 34044.668
$ws.Range("L123").Value = 34044.668
$ws.Range("N123").Value = -43844.668
$ws.Range("H135").Value = 53572464
$ws.Range("I135").Value = 38462564
$ws.Range("J135").Value = 250001150
$ws.Range("K135").Value = 346163076
$ws.Range("L135").Value = 2250010350
$ws.Range("M135").Value = -346160541
$ws.Range("N135").Value = -2250015420
$ws.Range("H141").Value = 6484
$ws.Range("I141").Value = 4700
$ws.Range("J141").Value = 10052
$ws.Range("K141").Value = 14100
$ws.Range("L141").Value = 30156
$ws.Range("M141").Value = -8920
$ws.Range("N141").Value = -40516

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H97").Value = 975.2692
$ws.Range("I97").Value = 774.9375
$ws.Range("J97").Value = 1295.8
$ws.Range("K97").Value = 774.9375
$ws.Range("L97").Value = 1295.8
$ws.Range("M97").Value = -278.9375
$ws.Range("N97").Value = -2287.8
$ws.Range("H102").Value = 14112.789
$ws.Range("I102").Value = 2415
$ws.Range("J102").Value = 19511.77
$ws.Range("K102").Value = 2415
$ws.Range("L102").Value = 19511.77
$ws.Range("M102").Value = -793
$ws.Range("N102").Value = -22755.77

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 566.1852
$ws.Range("I94").Value = 564.5
$ws.Range("J94").Value = 610
$ws.Range("K94").Value = 564.5
$ws.Range("L94").Value = 610
$ws.Range("M94").Value = -113.5
$ws.Range("N94").Value = -1512
$ws.Range("H99").Value = 1925.8636
$ws.Range("I99").Value = 1874.1904
$ws.Range("K99").Value = 1874.1904
$ws.Range("M99").Value = -376.1904
$ws.Range("H105").Value = 3328.1
$ws.Range("I105").Value = 2988.5715
$ws.Range("J105").Value = 4120.3335
$ws.Range("K105").Value = 2988.5715
$ws.Range("L105").Value = 4120.3335
$ws.Range("M105").Value = -1241.5715
$ws.Range("N105").Value = -7614.3335
$ws.Range("H134").Value = 3690.988
$ws.Range("I134").Value = 1932.8695
$ws.Range("K134").Value = 5798.6085
$ws.Range("M134").Value = -3263.6085

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1531.0625
$ws.Range("I22").Value = 222.4
$ws.Range("J22").Value = 2125.9092
$ws.Range("K22").Value = 222.4
$ws.Range("L22").Value = 2125.9092
$ws.Range("M22").Value = 127.6
$ws.Range("N22").Value = -2825.9092
$ws.Range("H31").Value = 171156.97
$ws.Range("I31").Value = 2194.2
$ws.Range("J31").Value = 207888.02
$ws.Range("K31").Value = 2194.2
$ws.Range("L31").Value = 207888.02
$ws.Range("M31").Value = -1899.2
$ws.Range("N31").Value = -208478.02
$ws.Range("H34").Value = 171156.97
$ws.Range("I34").Value = 2194.2
$ws.Range("J34").Value = 207888.02
$ws.Range("K34").Value = 2194.2
$ws.Range("L34").Value = 207888.02
$ws.Range("M34").Value = -1992.2
$ws.Range("N34").Value = -208292.02
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H114").Value = 23421
$ws.Range("J114").Value = 23421
$ws.Range("L114").Value = 23421
$ws.Range("N114").Value = -32099
$ws.Range("H120").Value = 30756
$ws.Range("J120").Value = 30756
$ws.Range("L120").Value = 30756
$ws.Range("N120").Value = -38014
$ws.Range("H130").Value = 39374
$ws.Range("J130").Value = 39374
$ws.Range("L130").Value = 39374
$ws.Range("N130").Value = -49414
$ws.Range("H134").Value = 1169369.6
$ws.Range("I134").Value = 1552.75
$ws.Range("J134").Value = 1753278.1
$ws.Range("K134").Value = 4658.25
$ws.Range("L134").Value = 5259834.300000001
$ws.Range("M134").Value = -2123.25
$ws.Range("N134").Value = -5264904.300000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 9564.666999999999
$ws.Range("I122").Value = 643
$ws.Range("J122").Value = 14025.5
$ws.Range("K122").Value = 5787
$ws.Range("L122").Value = 126229.5
$ws.Range("M122").Value = -3337
$ws.Range("N122").Value = -131129.5
$ws.Range("H132").Value = 623.75
$ws.Range("J132").Value = 497.5
$ws.Range("L132").Value = 4477.5
$ws.Range("N132").Value = -9537.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H55").Value = 8000
$ws.Range("J55").Value = 8000
$ws.Range("L55").Value = 8000
$ws.Range("N55").Value = -8654
$ws.Range("H97").Value = 1678.9445
$ws.Range("I97").Value = 1333
$ws.Range("J97").Value = 2111.375
$ws.Range("K97").Value = 1333
$ws.Range("L97").Value = 2111.375
$ws.Range("M97").Value = -837
$ws.Range("N97").Value = -3103.375
$ws.Range("H132").Value = 4087.6316
$ws.Range("I132").Value = 2008.7778
$ws.Range("K132").Value = 6026.3334
$ws.Range("M132").Value = -3496.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2899.158
$ws.Range("J40").Value = 4080.8333
$ws.Range("L40").Value = 4080.8333
$ws.Range("N40").Value = -4352.8333
$ws.Range("H45").Value = 11895
$ws.Range("I45").Value = 10000
$ws.Range("J45").Value = 13790
$ws.Range("K45").Value = 10000
$ws.Range("L45").Value = 13790
$ws.Range("M45").Value = -9593
$ws.Range("N45").Value = -14604
$ws.Range("H93").Value = 1352.5927
$ws.Range("I93").Value = 1200.4546
$ws.Range("J93").Value = 1457.1875
$ws.Range("K93").Value = 1200.4546
$ws.Range("L93").Value = 1457.1875
$ws.Range("M93").Value = 47.54539999999997
$ws.Range("N93").Value = -3953.1875
$ws.Range("H122").Value = 2072.8
$ws.Range("I122").Value = 2078
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 6234
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -3784
$ws.Range("N122").Value = -10900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 950
$ws.Range("J14").Value = 950
$ws.Range("L14").Value = 950
$ws.Range("N14").Value = -1286
$ws.Range("H19").Value = 83339.336
$ws.Range("J19").Value = 83339.336
$ws.Range("L19").Value = 83339.336
$ws.Range("N19").Value = -83687.336
$ws.Range("H96").Value = 1427
$ws.Range("I96").Value = 832.3333
$ws.Range("J96").Value = 1681.8572
$ws.Range("K96").Value = 832.3333
$ws.Range("L96").Value = 1681.8572
$ws.Range("M96").Value = 540.6667
$ws.Range("N96").Value = -4427.8572
$ws.Range("H100").Value = 470.36365
$ws.Range("I100").Value = 419.1
$ws.Range("J100").Value = 983
$ws.Range("K100").Value = 838.2
$ws.Range("L100").Value = 1966
$ws.Range("M100").Value = -297.2
$ws.Range("N100").Value = -3048
$ws.Range("H114").Value = 14296.8
$ws.Range("J114").Value = 14296.8
$ws.Range("L114").Value = 14296.8
$ws.Range("N114").Value = -22974.8
